$d = $word.ActiveDocument

# --- 1. Insert a new "Meta description" paragraph right after the H1 title ---
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:r/>' +
           '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
           '<w:r><w:t>: Get the nostalgic feel of traditional slot symbols with Diamond 7 by Novomatic. Play now and multiply your bet up to 2500 times, for free!</w:t></w:r>' +
           '</w:p>'
$metaPara.Range.InsertXML($metaXml)

# --- 2. Near the end: drop the duplicated bold title paragraph, and replace ---
#        the final (italic) paragraph's text with the new image prompt.
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitlePara.Range.Delete() | Out-Null

$count2 = $d.Paragraphs.Count
$imgPromptPara = $d.Paragraphs.Item($count2)
$imgXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:r/>' +
          '<w:r><w:rPr><w:i/></w:rPr><w:t>Create a cartoon-style feature image for &quot;Diamond 7&quot; that features a happy Maya warrior with glasses. The warrior should be holding diamonds and surrounded by traditional slot symbols like the number 7, cherries, and the BAR sign. The image should be bright and colorful, with a fun and playful vibe that reflects the simplicity and straightforwardness of the game. The overall style should be cartoonish and eye-catching to draw in potential players who enjoy classic slot games.</w:t></w:r>' +
          '</w:p>'
$imgPromptPara.Range.InsertXML($imgXml)
